# [Improvement] On terminology : room -> bed
$wb = $excel.ActiveWorkbook

# --- Rename the "rooms" sheet to "beds" ---
$wsBabies = $wb.Worksheets.Item("babies")
$wsBeds   = $wb.Worksheets.Item("rooms")
$wsBeds.Name = "beds"

# --- Update room -> bed terminology on the (now) "beds" sheet header row ---
$wsBeds.Range("A1").Value = "all_beds"
$wsBeds.Range("B1").Value = "new_beds"
$wsBeds.Range("C1").Value = "old_beds"
$wsBeds.Range("E1").Value = "new_beds_service"
$wsBeds.Range("F1").Value = "old_beds_service"
$wsBeds.Range("G1").Value = "beds_capacities"

# --- Update selections / active sheet to match the new authoring state ---
$wsBabies.Range("C13").Select()
$wsBeds.Range("E16").Select()
